$p = $ppt.ActivePresentation

# The standalone "Picture" slide (currently at position 11) is moved to
# position 8, i.e. right after the "Class Diagram - Chess Pieces" slide and
# before the "Library Functions Used" / "User-Defined Functions" slides.
$p.Slides.Item(11).MoveTo(8)
